$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("usermaint")

# Update the username value in row 2 (was "jsmith", now "jsmith2")
$ws.Range("C2").Value = "jsmith2"

# Update the active selection to reflect the edited cell
$ws.Range("C2").Select()
